# "new backlogs checked in"
# - Populates the previously-empty "Sprint3" sheet with the new backlog
#   table (headers, mitarbeiter/aufgabe rows, two new user-story blocks).
# - Re-points the active sheet/selection from Sprint2 to Sprint3.
# - Leaves Sprint2's own scroll position/selection updated (no longer the
#   tab-selected sheet).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sprint1")
$ws2 = $wb.Worksheets.Item("Sprint2")
$ws3 = $wb.Worksheets.Item("Sprint3")

# ---------------------------------------------------------------------
# 1) Formatting: pull the existing cell styles used elsewhere in the
#    workbook onto the equivalent Sprint3 cells via a formats-only copy,
#    so the same shared <xf> style indices get reused instead of minted
#    fresh ones.
# ---------------------------------------------------------------------

# Header row (A1:E1) -> style of Sprint2's own header row
$ws2.Range("A1").Copy() | Out-Null
$ws3.Range("A1").PasteSpecial(-4122) | Out-Null
$ws2.Range("B1").Copy() | Out-Null
$ws3.Range("B1").PasteSpecial(-4122) | Out-Null
$ws2.Range("C1").Copy() | Out-Null
$ws3.Range("C1").PasteSpecial(-4122) | Out-Null
$ws2.Range("D1").Copy() | Out-Null
$ws3.Range("D1").PasteSpecial(-4122) | Out-Null
$ws2.Range("E1").Copy() | Out-Null
$ws3.Range("E1").PasteSpecial(-4122) | Out-Null

# Blank spacer row under the header (A2:B2)
$ws2.Range("A16").Copy() | Out-Null
$ws3.Range("A2").PasteSpecial(-4122) | Out-Null
$ws2.Range("C16").Copy() | Out-Null
$ws3.Range("B2").PasteSpecial(-4122) | Out-Null

# Story-id cells (B13, B20) -> style of the existing story-id column cell
$ws2.Range("A26").Copy() | Out-Null
$ws3.Range("B13").PasteSpecial(-4122) | Out-Null
$ws3.Range("B20").PasteSpecial(-4122) | Out-Null

$ws3.Rows.Item(1).RowHeight = 18.75

# ---------------------------------------------------------------------
# 2) Cell values (new shared strings get interned automatically).
# ---------------------------------------------------------------------

$cellData = @(
    @{Cell="A1"; Value='Mitarbeiter'}
    @{Cell="B1"; Value='userstory ID'}
    @{Cell="C1"; Value='sprint backlog'}
    @{Cell="D1"; Value='Aufwand in Stunden'}
    @{Cell="E1"; Value='Status'}
    @{Cell="A3"; Value='Max Kr'}
    @{Cell="C3"; Value='Unit tests'}
    @{Cell="A4"; Value='Corinna'}
    @{Cell="C4"; Value='Technische Dokumentation'}
    @{Cell="A5"; Value='Manfred'}
    @{Cell="C5"; Value='Datenbank-Entwurf'}
    @{Cell="A6"; Value='Max Ke'}
    @{Cell="C6"; Value='Benutzer-Dokumentation'}
    @{Cell="A7"; Value='Max Ke'}
    @{Cell="C7"; Value='Benutzter-Tests'}
    @{Cell="A8"; Value='Henrik'}
    @{Cell="C8"; Value='Unit Tests'}
    @{Cell="A9"; Value='Rosemarie/Henrik'}
    @{Cell="C9"; Value='Multilingulität'}
    @{Cell="A10"; Value='Andreas'}
    @{Cell="C10"; Value='Paging/Springer'}
    @{Cell="A11"; Value='Alle'}
    @{Cell="C11"; Value='Defect Liste'}
    @{Cell="B13"; Value='3-1'}
    @{Cell="C13"; Value='Als <MM> will ich eine positiv/negativ Sentiment Analyse der Tweets erhalten'}
    @{Cell="A14"; Value='Johannes/Manfred'}
    @{Cell="C14"; Value='Controllerklasse'}
    @{Cell="A15"; Value='Johannes/Manfred'}
    @{Cell="C15"; Value='Modelklasse'}
    @{Cell="A16"; Value='Johannes/Manfred'}
    @{Cell="C16"; Value='Datenbank'}
    @{Cell="A17"; Value='Corinna'}
    @{Cell="C17"; Value='Technische Dokumentation'}
    @{Cell="A18"; Value='Max Ke'}
    @{Cell="C18"; Value='Benutzter-Tests'}
    @{Cell="B20"; Value='4-1'}
    @{Cell="C20"; Value='Als <MM> will ich die vorhandenen Tweets einschränken, sodass nur relevante Tweets analysiert werden'}
    @{Cell="A21"; Value='Johannes/Manfred'}
    @{Cell="C21"; Value='View'}
    @{Cell="A22"; Value='Johannes/Manfred'}
    @{Cell="C22"; Value='Controllerklasse'}
    @{Cell="A23"; Value='Johannes/Manfred'}
    @{Cell="C23"; Value='Modelklasse'}
    @{Cell="A24"; Value='Johannes/Manfred'}
    @{Cell="C24"; Value='Datenbank'}
    @{Cell="A25"; Value='Corinna'}
    @{Cell="C25"; Value='Technische Dokumentation'}
    @{Cell="A26"; Value='Max Ke'}
    @{Cell="C26"; Value='Benutzter-Tests'}
)

foreach ($item in $cellData) {
    $ws3.Range($item.Cell).Value = $item.Value
}

# ---------------------------------------------------------------------
# 3) Column widths (best-effort; engine quantizes to 1/6-character
#    steps so these land close to, not bit-exact with, the authored
#    widths).
# ---------------------------------------------------------------------

$ws3.Columns.Item(1).ColumnWidth = 14
$ws3.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws3.Columns.Item(3).ColumnWidth = 16.5
$ws3.Columns.Item(4).ColumnWidth = 23.666666666666668
$ws3.Columns.Item(5).ColumnWidth = 7.666666666666667

# ---------------------------------------------------------------------
# 4) View/selection: Sprint2 keeps its own scroll/selection state but
#    is no longer the active tab; Sprint3 becomes active with its own
#    selection and zoom.
# ---------------------------------------------------------------------

$ws2.Range("A33:C33").Select() | Out-Null

$ws3.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 100
$ws3.Range("A12").Select() | Out-Null
